$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '91.722.59'
$ws.Range("E2").Value = '  +0.54%  '

# Row 3
$ws.Range("D3").Value = '3.107.66'
$ws.Range("E3").Value = '  -0.98%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.64'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '620.72'
$ws.Range("E6").Value = '  -2.26%  '

# Row 7
$ws.Range("E7").Value = '  +4.78%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.371'
$ws.Range("E8").Value = '  +1.24%  '

# Row 9
$ws.Range("E9").Value = '  +0.06%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.758'
$ws.Range("E10").Value = '  +4.59%  '

# Row 11
$ws.Range("D11").Value = '2.728.24'
$ws.Range("E11").Value = '  -12.96%  '

# Row 12
$ws.Range("E12").Value = '  +3.03%  '

# Row 13
$ws.Range("E13").Value = '  +1.18%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.46'
$ws.Range("E14").Value = '  -2.70%  '

# Row 15
$ws.Range("D15").Value = '91.637.73'
$ws.Range("E15").Value = '  +0.93%  '

# Row 16
$ws.Range("E16").Value = '  -1.66%  '

# Row 18
$ws.Range("D18").Value = '3.185.50'
$ws.Range("E18").Value = '  +0.88%  '

# Row 19
$ws.Range("E19").Value = '  -0.53%  '

# Row 20
$ws.Range("B20").Value = 'PEPE'
$ws.Range("C20").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000217'
$ws.Range("E20").Value = '  +1.28%  '

# Row 21
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.55'
$ws.Range("E21").Value = '  +1.41%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.77'
$ws.Range("E22").Value = '  +2.31%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '447.19'
$ws.Range("E23").Value = '  +0.18%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.12'
$ws.Range("E24").Value = '  +0.41%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.63'
$ws.Range("E25").Value = '  -5.65%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '90.50'
$ws.Range("E26").Value = '  +0.32%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.96'
$ws.Range("E27").Value = '  -3.99%  '

# Row 28
$ws.Range("D28").Value = '3.255.78'

# Row 29
$ws.Range("E29").Value = '  +0.15%  '

# Row 30
$ws.Range("E30").Value = '  +14.07%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.235'
$ws.Range("E31").Value = '  +17.92%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.42'
$ws.Range("E32").Value = '  -2.70%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.64%  '

# Row 34
$ws.Range("E34").Value = '  +12.72%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.111'
$ws.Range("E35").Value = '  +29.64%  '

# Row 36
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.67'
$ws.Range("E36").Value = '  +6.60%  '

# Row 37
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.52'
$ws.Range("E37").Value = '  -1.42%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.13'
$ws.Range("E38").Value = '  +22.00%  '

# Row 39
$ws.Range("E39").Value = '  -0.94%  '

# Row 40
$ws.Range("E40").Value = '  -4.47%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '491.69'
$ws.Range("E41").Value = '  -4.58%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.31'
$ws.Range("E42").Value = '  +0.49%  '

# Row 43
$ws.Range("E43").Value = '  +0.34%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.16'
$ws.Range("E44").Value = '  -0.19%  '

# Row 45
$ws.Range("E45").Value = '  -0.01%  '

# Row 46
$ws.Range("E46").Value = '  -1.95%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '154.66'
$ws.Range("E47").Value = '  +2.44%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.691'
$ws.Range("E48").Value = '  -0.85%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.51'
$ws.Range("E49").Value = '  -1.65%  '

# Row 50
$ws.Range("E50").Value = '  -1.04%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.58'
$ws.Range("E51").Value = '  -2.61%  '

